$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.238.41'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.68%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.513.63'
$ws.Range('D3').Style = 'Normal'

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.08%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '572.73'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.47%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.52'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.60%  '

$ws.Range('E7').Value = '  +0.02%  '

$ws.Range('E8').Value = '  +1.54%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.511.80'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.49%  '

$ws.Range('E10').Value = '  -0.76%  '

$ws.Range('E11').Value = '  -0.41%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.354'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +3.52%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.92'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.56%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.975.33'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -0.35%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '69.159.14'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.55%  '

$ws.Range('E16').Value = '  -2.39%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.84'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.27%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.516.09'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.57%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.36'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.88%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.61'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.25%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '349.22'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.40%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '3.93'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.07%  '

$ws.Range('E23').Value = '  +1.27%  '

$ws.Range('E24').Value = '  -0.06%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '70.41'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.92%  '

$ws.Range('E26').Value = '  -1.38%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.96'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.28%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.650.06'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -0.10%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.998'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.32%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0891'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.75%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.85'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.08%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '462.96'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.02%  '

$ws.Range('E33').Value = '  -2.70%  '

$ws.Range('E34').Value = '  -1.71%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.999'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.05%  '

$ws.Range('E36').Value = '  +2.13%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '158.08'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.46%  '

$ws.Range('E38').Value = '  +0.85%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.52'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.10%  '

$ws.Range('E40').Value = '  -0.04%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.75'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.24%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.319'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.08%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.61'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.19%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '38.04'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.62%  '

$ws.Range('E45').Value = '  -13.35%  '

$ws.Range('E46').Value = '  -5.26%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '141.79'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.34%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.526'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.79%  '

$ws.Range('E49').Value = '  -1.40%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0730'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.37%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.56'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.88%  '
